$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Fitness values (column C) for the corresponding rows as per the
# updated run_9.xlsx log data.
$ws.Range("C2").Value = 13256
$ws.Range("C3").Value = 12017
$ws.Range("C4:C6").Value = 10780
$ws.Range("C7:C10").Value = 10352
$ws.Range("C11:C16").Value = 10139
$ws.Range("C17").Value = 10082
$ws.Range("C18:C22").Value = 9994
$ws.Range("C23:C32").Value = 9330
$ws.Range("C33").Value = 9207
$ws.Range("C34:C37").Value = 9176
$ws.Range("C38:C40").Value = 7917
$ws.Range("C41:C42").Value = 7815
$ws.Range("C43:C52").Value = 7734
$ws.Range("C84:C90").Value = 7734
$ws.Range("C91:C92").Value = 7685
$ws.Range("C113:C134").Value = 7343
